$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C) column date serial from 45208 to 45212 for rows 2-99
for ($r = 2; $r -le 99; $r++) {
    $ws.Cells.Item($r, 3).Value = 45212
}

# Update hyperlink formulas (S,T,V,W,X,Y) for rows 2-5 to point at renamed files

# Row 2
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/artfynd/A 30683-2023 artfynd.xlsx", "A 30683-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/kartor/A 30683-2023 karta.png", "A 30683-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomål/A 30683-2023 fsc-klagomål.docx", "A 30683-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomålsmail/A 30683-2023 fsc-klagomål mail.docx", "A 30683-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/tillsyn/A 30683-2023 tillsynsbegäran.docx", "A 30683-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/ti,llsynsmail/A 30683-2023 tillsynsbegäran mail.docx", "A 30683-2023")'

# Row 3
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/artfynd/A 32699-2023 artfynd.xlsx", "A 32699-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/kartor/A 32699-2023 karta.png", "A 32699-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomål/A 32699-2023 fsc-klagomål.docx", "A 32699-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomålsmail/A 32699-2023 fsc-klagomål mail.docx", "A 32699-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/tillsyn/A 32699-2023 tillsynsbegäran.docx", "A 32699-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/ti,llsynsmail/A 32699-2023 tillsynsbegäran mail.docx", "A 32699-2023")'

# Row 4
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/artfynd/A 29992-2023 artfynd.xlsx", "A 29992-2023")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/kartor/A 29992-2023 karta.png", "A 29992-2023")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomål/A 29992-2023 fsc-klagomål.docx", "A 29992-2023")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/klagomålsmail/A 29992-2023 fsc-klagomål mail.docx", "A 29992-2023")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/tillsyn/A 29992-2023 tillsynsbegäran.docx", "A 29992-2023")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2380/ti,llsynsmail/A 29992-2023 tillsynsbegäran mail.docx", "A 29992-2023")'

# Row 5
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/artfynd/A 30834-2023 artfynd.xlsx", "A 30834-2023")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/kartor/A 30834-2023 karta.png", "A 30834-2023")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/klagomål/A 30834-2023 fsc-klagomål.docx", "A 30834-2023")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/klagomålsmail/A 30834-2023 fsc-klagomål mail.docx", "A 30834-2023")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/tillsyn/A 30834-2023 tillsynsbegäran.docx", "A 30834-2023")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/ti,llsynsmail/A 30834-2023 tillsynsbegäran mail.docx", "A 30834-2023")'
